$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-18 12:49:36"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
